$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C60").Value = 7310
$ws.Range("C61:C72").Value = 7295
$ws.Range("C82:C109").Value = 7293
